# This script rewrites the monthly recurrence/retention metrics table (A2:D82)
# with the refreshed dataset (now spanning 2018-08 through 2025-04).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = New-Object "object[,]" 81,4
$data[0,0] = "2018-08"
$data[0,1] = 2
$data[0,2] = 156
$data[0,3] = 1.282051282051282
$data[1,0] = "2018-09"
$data[1,1] = 4
$data[1,2] = 195
$data[1,3] = 2.051282051282051
$data[2,0] = "2018-10"
$data[2,1] = 5
$data[2,2] = 177
$data[2,3] = 2.824858757062147
$data[3,0] = "2018-11"
$data[3,1] = 7
$data[3,2] = 187
$data[3,3] = 3.74331550802139
$data[4,0] = "2018-12"
$data[4,1] = 2
$data[4,2] = 170
$data[4,3] = 1.176470588235294
$data[5,0] = "2019-01"
$data[5,1] = 3
$data[5,2] = 179
$data[5,3] = 1.675977653631285
$data[6,0] = "2019-02"
$data[6,1] = 1
$data[6,2] = 165
$data[6,3] = 0.6060606060606061
$data[7,0] = "2019-03"
$data[7,1] = 7
$data[7,2] = 136
$data[7,3] = 5.147058823529411
$data[8,0] = "2019-04"
$data[8,1] = 5
$data[8,2] = 121
$data[8,3] = 4.132231404958678
$data[9,0] = "2019-05"
$data[9,1] = 3
$data[9,2] = 112
$data[9,3] = 2.678571428571428
$data[10,0] = "2019-06"
$data[10,1] = 2
$data[10,2] = 99
$data[10,3] = 2.02020202020202
$data[11,0] = "2019-07"
$data[11,1] = 3
$data[11,2] = 117
$data[11,3] = 2.564102564102564
$data[12,0] = "2019-08"
$data[12,1] = 7
$data[12,2] = 127
$data[12,3] = 5.511811023622047
$data[13,0] = "2019-09"
$data[13,1] = 11
$data[13,2] = 232
$data[13,3] = 4.741379310344827
$data[14,0] = "2019-10"
$data[14,1] = 3
$data[14,2] = 227
$data[14,3] = 1.3215859030837
$data[15,0] = "2019-11"
$data[15,1] = 8
$data[15,2] = 381
$data[15,3] = 2.099737532808399
$data[16,0] = "2019-12"
$data[16,1] = 8
$data[16,2] = 330
$data[16,3] = 2.424242424242424
$data[17,0] = "2020-01"
$data[17,1] = 10
$data[17,2] = 344
$data[17,3] = 2.906976744186046
$data[18,0] = "2020-02"
$data[18,1] = 8
$data[18,2] = 301
$data[18,3] = 2.6578073089701
$data[19,0] = "2020-03"
$data[19,1] = 4
$data[19,2] = 235
$data[19,3] = 1.702127659574468
$data[20,0] = "2020-04"
$data[20,1] = 8
$data[20,2] = 207
$data[20,3] = 3.864734299516908
$data[21,0] = "2020-05"
$data[21,1] = 4
$data[21,2] = 161
$data[21,3] = 2.484472049689441
$data[22,0] = "2020-06"
$data[22,1] = 5
$data[22,2] = 234
$data[22,3] = 2.136752136752137
$data[23,0] = "2020-07"
$data[23,1] = 2
$data[23,2] = 272
$data[23,3] = 0.7352941176470588
$data[24,0] = "2020-08"
$data[24,1] = 8
$data[24,2] = 294
$data[24,3] = 2.72108843537415
$data[25,0] = "2020-09"
$data[25,1] = 10
$data[25,2] = 317
$data[25,3] = 3.154574132492113
$data[26,0] = "2020-10"
$data[26,1] = 10
$data[26,2] = 307
$data[26,3] = 3.257328990228013
$data[27,0] = "2020-11"
$data[27,1] = 11
$data[27,2] = 254
$data[27,3] = 4.330708661417323
$data[28,0] = "2020-12"
$data[28,1] = 7
$data[28,2] = 230
$data[28,3] = 3.043478260869565
$data[29,0] = "2021-01"
$data[29,1] = 2
$data[29,2] = 213
$data[29,3] = 0.9389671361502347
$data[30,0] = "2021-02"
$data[30,1] = 1
$data[30,2] = 121
$data[30,3] = 0.8264462809917356
$data[31,0] = "2021-03"
$data[31,1] = 2
$data[31,2] = 157
$data[31,3] = 1.273885350318471
$data[32,0] = "2021-04"
$data[32,1] = 6
$data[32,2] = 253
$data[32,3] = 2.371541501976284
$data[33,0] = "2021-05"
$data[33,1] = 7
$data[33,2] = 236
$data[33,3] = 2.966101694915254
$data[34,0] = "2021-06"
$data[34,1] = 12
$data[34,2] = 248
$data[34,3] = 4.838709677419355
$data[35,0] = "2021-07"
$data[35,1] = 12
$data[35,2] = 252
$data[35,3] = 4.761904761904762
$data[36,0] = "2021-08"
$data[36,1] = 9
$data[36,2] = 280
$data[36,3] = 3.214285714285714
$data[37,0] = "2021-09"
$data[37,1] = 13
$data[37,2] = 230
$data[37,3] = 5.652173913043478
$data[38,0] = "2021-10"
$data[38,1] = 7
$data[38,2] = 223
$data[38,3] = 3.139013452914798
$data[39,0] = "2021-11"
$data[39,1] = 6
$data[39,2] = 218
$data[39,3] = 2.752293577981652
$data[40,0] = "2021-12"
$data[40,1] = 9
$data[40,2] = 292
$data[40,3] = 3.082191780821918
$data[41,0] = "2022-01"
$data[41,1] = 12
$data[41,2] = 256
$data[41,3] = 4.6875
$data[42,0] = "2022-02"
$data[42,1] = 7
$data[42,2] = 254
$data[42,3] = 2.755905511811024
$data[43,0] = "2022-03"
$data[43,1] = 9
$data[43,2] = 236
$data[43,3] = 3.813559322033898
$data[44,0] = "2022-04"
$data[44,1] = 14
$data[44,2] = 269
$data[44,3] = 5.204460966542751
$data[45,0] = "2022-05"
$data[45,1] = 11
$data[45,2] = 244
$data[45,3] = 4.508196721311475
$data[46,0] = "2022-06"
$data[46,1] = 12
$data[46,2] = 282
$data[46,3] = 4.25531914893617
$data[47,0] = "2022-07"
$data[47,1] = 13
$data[47,2] = 261
$data[47,3] = 4.980842911877394
$data[48,0] = "2022-08"
$data[48,1] = 14
$data[48,2] = 279
$data[48,3] = 5.017921146953405
$data[49,0] = "2022-09"
$data[49,1] = 23
$data[49,2] = 301
$data[49,3] = 7.641196013289036
$data[50,0] = "2022-10"
$data[50,1] = 11
$data[50,2] = 293
$data[50,3] = 3.754266211604096
$data[51,0] = "2022-11"
$data[51,1] = 11
$data[51,2] = 266
$data[51,3] = 4.135338345864661
$data[52,0] = "2022-12"
$data[52,1] = 11
$data[52,2] = 237
$data[52,3] = 4.641350210970464
$data[53,0] = "2023-01"
$data[53,1] = 6
$data[53,2] = 253
$data[53,3] = 2.371541501976284
$data[54,0] = "2023-02"
$data[54,1] = 15
$data[54,2] = 296
$data[54,3] = 5.067567567567568
$data[55,0] = "2023-03"
$data[55,1] = 7
$data[55,2] = 264
$data[55,3] = 2.651515151515151
$data[56,0] = "2023-04"
$data[56,1] = 8
$data[56,2] = 271
$data[56,3] = 2.952029520295203
$data[57,0] = "2023-05"
$data[57,1] = 4
$data[57,2] = 233
$data[57,3] = 1.716738197424893
$data[58,0] = "2023-06"
$data[58,1] = 15
$data[58,2] = 293
$data[58,3] = 5.119453924914676
$data[59,0] = "2023-07"
$data[59,1] = 10
$data[59,2] = 263
$data[59,3] = 3.802281368821293
$data[60,0] = "2023-08"
$data[60,1] = 11
$data[60,2] = 249
$data[60,3] = 4.417670682730924
$data[61,0] = "2023-09"
$data[61,1] = 8
$data[61,2] = 240
$data[61,3] = 3.333333333333333
$data[62,0] = "2023-10"
$data[62,1] = 6
$data[62,2] = 285
$data[62,3] = 2.105263157894737
$data[63,0] = "2023-11"
$data[63,1] = 5
$data[63,2] = 209
$data[63,3] = 2.392344497607656
$data[64,0] = "2023-12"
$data[64,1] = 9
$data[64,2] = 224
$data[64,3] = 4.017857142857143
$data[65,0] = "2024-01"
$data[65,1] = 8
$data[65,2] = 259
$data[65,3] = 3.088803088803089
$data[66,0] = "2024-02"
$data[66,1] = 13
$data[66,2] = 232
$data[66,3] = 5.603448275862069
$data[67,0] = "2024-03"
$data[67,1] = 8
$data[67,2] = 210
$data[67,3] = 3.80952380952381
$data[68,0] = "2024-04"
$data[68,1] = 17
$data[68,2] = 237
$data[68,3] = 7.172995780590717
$data[69,0] = "2024-05"
$data[69,1] = 20
$data[69,2] = 241
$data[69,3] = 8.298755186721991
$data[70,0] = "2024-06"
$data[70,1] = 16
$data[70,2] = 191
$data[70,3] = 8.37696335078534
$data[71,0] = "2024-07"
$data[71,1] = 17
$data[71,2] = 207
$data[71,3] = 8.212560386473431
$data[72,0] = "2024-08"
$data[72,1] = 24
$data[72,2] = 267
$data[72,3] = 8.988764044943821
$data[73,0] = "2024-09"
$data[73,1] = 16
$data[73,2] = 263
$data[73,3] = 6.083650190114068
$data[74,0] = "2024-10"
$data[74,1] = 12
$data[74,2] = 230
$data[74,3] = 5.217391304347826
$data[75,0] = "2024-11"
$data[75,1] = 18
$data[75,2] = 257
$data[75,3] = 7.003891050583658
$data[76,0] = "2024-12"
$data[76,1] = 17
$data[76,2] = 264
$data[76,3] = 6.439393939393939
$data[77,0] = "2025-01"
$data[77,1] = 16
$data[77,2] = 257
$data[77,3] = 6.22568093385214
$data[78,0] = "2025-02"
$data[78,1] = 15
$data[78,2] = 230
$data[78,3] = 6.521739130434782
$data[79,0] = "2025-03"
$data[79,1] = 14
$data[79,2] = 225
$data[79,3] = 6.222222222222222
$data[80,0] = "2025-04"
$data[80,1] = 12
$data[80,2] = 228
$data[80,3] = 5.263157894736842

$ws.Range("A2:D82").Value = $data

